$d = $word.ActiveDocument

# 1. Merge the split "Cuando se llena..." sentence back into a single run.
#    This Find/Replace spans the old _GoBack bookmark's position, which
#    removes that bookmark as a side effect (matches the diff: the
#    bookmark is gone from its old spot, the two runs become one).
$d.Content.Find.Execute(
    "Cuando se llena el formulario con toda la información y dar en el botón guardar, la información no está quedando almacenada en la base de datos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Cuando se llena el formulario con toda la información y dar en el botón guardar, la información no está quedando almacenada en la base de datos",
    2)

# 2. All three "Pendiente" statuses become "Resuelta".
$find = $d.Content.Find
$find.ClearFormatting()
while ($find.Execute("Pendiente", $true, $false, $false, $false, $false, $true, 1, $false, "Resuelta", 2)) {
}

# 3. Re-insert the _GoBack bookmark right after the text of the last
#    "Resuelta" (the status that used to be the final "Pendiente").
#    Locate the end of that last occurrence.
$search = $d.Content
$search.Find.ClearFormatting()
$lastEnd = -1
while ($search.Find.Execute("Resuelta", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastEnd = $search.End
    $search.Collapse(0)
}

# Placing a bookmark directly on a collapsed Range that sits exactly at a
# paragraph boundary is unreliable, so nudge it into place: insert a
# placeholder character, anchor the bookmark across it (non-collapsed
# ranges land correctly), then remove the placeholder so the bookmark
# collapses back down to the original boundary position.
$insPoint = $d.Range($lastEnd, $lastEnd)
$insPoint.InsertAfter("X")
$bmRange = $d.Range($lastEnd, $lastEnd + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($lastEnd, $lastEnd + 1)
$placeholder.Text = ""

Write-Output "done"
